# Generate Report for Handback
# - Marks the 897ca6f8 file's handback status as failed on the Overview sheet.
# - Records the handback-transform error detail on both locale sheets (zh-cn, de-de).
# - Widens the "Error Detail" column so the new message is readable.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: update status for the 897ca6f8 row (row 7) ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E7").Value = "Handback transform failed"
$overview.Range("F7").Value = "Handback transform failed"

# --- zh-cn sheet: status, error detail + column width ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C7").Value = "Handback transform failed"
$zhcn.Range("P7").Value = "Handback file name: 4z4ypkyc.krx is different with handoff file name: 897ca6f8-c65e-49d2-a072-bd80d51d400d.e67a7238f2e04c3686c8a6bc59d98e6096946590.zh-cn."
$zhcn.Columns.Item(16).ColumnWidth = 235 / 6

# --- de-de sheet: status, error detail + column width ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C7").Value = "Handback transform failed"
$dede.Range("P7").Value = "Handback file name: 4z4ypkyc.krx is different with handoff file name: 897ca6f8-c65e-49d2-a072-bd80d51d400d.e67a7238f2e04c3686c8a6bc59d98e6096946590.de-de."
$dede.Columns.Item(16).ColumnWidth = 235 / 6
